$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "登录" (Login) - insert a new step between the old rows 6 and 7:
#   "创建token，包含phone" in red font at G7, everything below shifts down 1.
# ---------------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item(2)
$wsLogin.Rows.Item(7).Insert()
$wsLogin.Range("G7").Value = "创建token，包含phone"
$wsLogin.Range("G7").Font.Color = 255
$wsLogin.PageSetup.PaperSize = 9
$wsLogin.PageSetup.Orientation = 1
[void]$wsLogin.Activate()
[void]$wsLogin.Range("G21").Select()

# ---------------------------------------------------------------------------
# Sheet "用户设置" (User settings) - password flow fleshed out, email flow
# added.
# ---------------------------------------------------------------------------
$wsSettings = $wb.Worksheets.Item(3)

# Tab labels gain a space before "Tab"
$wsSettings.Range("C11").Value = "用户信息 Tab"
$wsSettings.Range("C15").Value = "密码 Tab"

# New date stamp on the "用户信息 Tab" block, matching A4's date style
$wsSettings.Range("A4").Copy()
$wsSettings.Range("A11").PasteSpecial(-4122)
$wsSettings.Range("A11").Value = 43574

# The old "api - updatePhone" note at D18 is relocated far below (new D31);
# clear its old slot first.
$wsSettings.Range("D18").ClearContents()

# Flesh out the password-update flow (rows 17-19)
$wsSettings.Range("G17").Value = "核对旧密码"
$wsSettings.Range("G18").Value = "更新密码"
$wsSettings.Range("D19").Value = "成功，清空表单"

# New "邮箱 Tab" (email) flow block (rows 21-28)
$wsSettings.Range("C21").Value = "邮箱 Tab"
$wsSettings.Range("F21").Value = "controller - Users"
$wsSettings.Range("D22").Value = "api - updateEmail"
$wsSettings.Range("G22").Value = "method - update_email"
$wsSettings.Range("G23").Value = "token中取出phone"
$wsSettings.Range("G24").Value = "核对验证码"
$wsSettings.Range("G25").Value = "以phone检索user"
$wsSettings.Range("G26").Value = "核对email"
$wsSettings.Range("G27").Value = "更新email"
$wsSettings.Range("D28").Value = "成功，清空表单，更新组件data.email"

# Relocated "api - updatePhone" note
$wsSettings.Range("D31").Value = "api - updatePhone"

[void]$wsSettings.Activate()
[void]$wsSettings.Range("G27").Select()

Write-Output "done"
